# Before pushing to github
# Updates the "Key " legend sheet and rebuilds the "Overall" report sheet.

$wb = $excel.ActiveWorkbook
$wsKey = $wb.Worksheets.Item("Key ")
$wsOverall = $wb.Worksheets.Item("Overall")

# ---------------------------------------------------------------------
# "Key " sheet: the red / blue(theme) / purple legend rows stay put;
# only the right-hand (G) column notes and the selected cell change.
# ---------------------------------------------------------------------
$wsKey.Activate()

$wsKey.Range("G1").Value = "NOT ON THE PRODUCT"
$wsKey.Range("G2").Value = "Strap(60) + Intertier(Varies) + Cell(60) "
$wsKey.Range("G3").Value = "Voltage is read every day "

$wsKey.Range("A2").Select()

# ---------------------------------------------------------------------
# "Overall" sheet: rebuilt from scratch with three labelled sections
# (Impedance / Voltage / Both), each a coloured header row followed by
# a data row of Substation / Actual / Required numbers.
# ---------------------------------------------------------------------
$wsOverall.Activate()
$wsOverall.Cells.Clear()

$wsOverall.Range("A1").Value = "Substation "
$wsOverall.Range("B1").Value = "Actual"
$wsOverall.Range("C1").Value = "Required"

$wsOverall.Range("A2").Value = "Impedance"
$wsOverall.Range("A2").Interior.Color = 255        # FF0000 red (reuses existing "Key " fill)

$wsOverall.Range("A3").Value = "Malborn"
$wsOverall.Range("B3").Value = 120
$wsOverall.Range("C3").Value = 130

$wsOverall.Range("A4").Value = "Voltage "
$wsOverall.Range("A4").Interior.Color = 12611584   # 0070C0 blue (new fill)

$wsOverall.Range("A5").Value = "Lauderdale West #1"
$wsOverall.Range("B5").Value = 1700
$wsOverall.Range("C5").Value = 1800

$wsOverall.Range("A6").Value = "Both"
$wsOverall.Range("A6").Interior.Color = 10498160   # 7030A0 purple (reuses existing "Key " fill)

$wsOverall.Range("A7").Value = "Shubuta"
$wsOverall.Range("B7").Value = 1900
$wsOverall.Range("C7").Value = 1930

$wsOverall.Range("E7").Select()
